$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4.44639156119888
$ws.Range("E2").Value = 0.0761924489145183
$ws.Range("F2").Value = 8.96780899395303
$ws.Range("D3").Value = -1.86516394141668
$ws.Range("E3").Value = -13.2224611827265
$ws.Range("F3").Value = 12.1858145640159
$ws.Range("D4").Value = -6.68901794185942
$ws.Range("E4").Value = -15.7954745928014
$ws.Range("F4").Value = 4.27847252222769
$ws.Range("D5").Value = -1.41548574072862
$ws.Range("E5").Value = -18.50287152233
$ws.Range("F5").Value = 21.6222072025779
$ws.Range("D6").Value = -46.0478700391514
$ws.Range("E6").Value = -52.3867737756833
$ws.Range("F6").Value = -38.2439332124747
$ws.Range("D7").Value = 20.3184764543881
$ws.Range("E7").Value = 13.3678609917057
$ws.Range("F7").Value = 28.0738401819582
$ws.Range("D8").Value = 21.3665435504278
$ws.Range("E8").Value = 11.7030444434211
$ws.Range("F8").Value = 32.4320452360681
$ws.Range("D9").Value = 58.3085560463384
$ws.Range("E9").Value = 31.879499682191
$ws.Range("F9").Value = 92.7027863266992
$ws.Range("D10").Value = 56.0790369227965
$ws.Range("E10").Value = 29.1024353677409
$ws.Range("F10").Value = 89.4324283429678
$ws.Range("D11").Value = -16.7507244450744
$ws.Range("E11").Value = -25.5953670402394
$ws.Range("F11").Value = -6.91032298469087
$ws.Range("D12").Value = 3.46560385189091
$ws.Range("E12").Value = -0.788958373958267
$ws.Range("F12").Value = 7.99332981257991
$ws.Range("D13").Value = 19.8519255937185
$ws.Range("E13").Value = 6.96088992711155
$ws.Range("F13").Value = 34.8723485490505
$ws.Range("D14").Value = 121.366351155699
$ws.Range("E14").Value = 62.3007149282717
$ws.Range("F14").Value = 218.09263915608
$ws.Range("D15").Value = 200.722797698643
$ws.Range("E15").Value = 104.674543316441
$ws.Range("F15").Value = 412.555646794438
$ws.Range("D16").Value = 147.573628946911
$ws.Range("E16").Value = 87.0798799347481
$ws.Range("F16").Value = 233.069670404704
$ws.Range("D17").Value = 4.02986896349452
$ws.Range("E17").Value = 0.911993207208108
$ws.Range("F17").Value = 7.57834103758783
$ws.Range("D18").Value = 21.963526903103
$ws.Range("E18").Value = 12.8331585417013
$ws.Range("F18").Value = 31.3956220592159
$ws.Range("D19").Value = 76.12957440222
$ws.Range("E19").Value = 53.5823162808539
$ws.Range("F19").Value = 105.083739896449
$ws.Range("D20").Value = 67.0693481217225
$ws.Range("E20").Value = 31.6880029015119
$ws.Range("F20").Value = 119.144091419085
$ws.Range("D21").Value = -43.5184735099785
$ws.Range("E21").Value = -49.0181133405591
$ws.Range("F21").Value = -36.8618504623877
